$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A2 to the value that used to be in A3, keeping it as text (inline string)
$ws.Range("A2").Value = "'79086"
$ws.Range("A2").Style = "Normal"

# Remove row 3 entirely (was: A3 = 79086 numeric)
$ws.Rows.Item(3).Delete()
